$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Insert a new blank column at H, shifting old H..M to I..N.
#    (This also auto-copies per-row cell styles from the old column H
#    (now I) onto the new blank column H, matching the target s= pattern.)
$ws.Columns("H").Insert()

# 2) Fill in the new ICON column header + type rows.
$ws.Range("H1").Value = "ICON"
$ws.Range("H2").Value = "l"
$ws.Range("H3").Value = "icon"
$ws.Range("H4").Value = "string"

# 3) Fill in the per-aura icon path values (5 rows per aura id).
$ws.Range("H5").Value = "icon/aura/aura_2001.png"
$ws.Range("H6").Value = "icon/aura/aura_2001.png"
$ws.Range("H7").Value = "icon/aura/aura_2001.png"
$ws.Range("H8").Value = "icon/aura/aura_2001.png"
$ws.Range("H9").Value = "icon/aura/aura_2001.png"

$ws.Range("H10").Value = "icon/aura/aura_2002.png"
$ws.Range("H11").Value = "icon/aura/aura_2002.png"
$ws.Range("H12").Value = "icon/aura/aura_2002.png"
$ws.Range("H13").Value = "icon/aura/aura_2002.png"
$ws.Range("H14").Value = "icon/aura/aura_2002.png"

$ws.Range("H15").Value = "icon/aura/aura_2003.png"
$ws.Range("H16").Value = "icon/aura/aura_2003.png"
$ws.Range("H17").Value = "icon/aura/aura_2003.png"
$ws.Range("H18").Value = "icon/aura/aura_2003.png"
$ws.Range("H19").Value = "icon/aura/aura_2003.png"

$ws.Range("H20").Value = "icon/aura/aura_2004.png"
$ws.Range("H21").Value = "icon/aura/aura_2004.png"
$ws.Range("H22").Value = "icon/aura/aura_2004.png"
$ws.Range("H23").Value = "icon/aura/aura_2004.png"
$ws.Range("H24").Value = "icon/aura/aura_2004.png"

$ws.Range("H25").Value = "icon/aura/aura_2005.png"
$ws.Range("H26").Value = "icon/aura/aura_2005.png"
$ws.Range("H27").Value = "icon/aura/aura_2005.png"
$ws.Range("H28").Value = "icon/aura/aura_2005.png"
$ws.Range("H29").Value = "icon/aura/aura_2005.png"

# 4) Widen the new ICON column. (Column I keeps the old column-H width of
#    10.375 automatically - the insert carries the <col> definition over.)
#    ColumnWidth is stored internally in 1/7-character increments, so the
#    input below is tuned to land on the closest representable value to
#    the target stored width of 27.25.
$ws.Columns("H").ColumnWidth = 26.55

# 5) The "script id" comment that used to sit on H1 stays pinned to the
#    H1 reference across the column insert (Excel does not follow moved
#    content), but its content now lives at I1 - so re-home the comment.
$oldComment = $ws.Range("H1").Comment
$commentText = $oldComment.Text()
$oldComment.Delete()
$ws.Range("I1").AddComment($commentText)

# 6) Match the saved selection/active cell from the target workbook.
$ws.Range("H32").Select()
